$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.404.43'
$ws.Range("E2").Value = '  -1.35%  '
$ws.Range("D3").Value = '1.710.12'
$ws.Range("E3").Value = '  -1.24%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.61'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5329'
$ws.Range("E6").Value = '  -1.95%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2670'
$ws.Range("E8").Value = '  -2.22%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06628'
$ws.Range("E9").Value = '  -0.69%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.93'
$ws.Range("E10").Value = '  -4.15%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07621'
$ws.Range("E11").Value = '  -2.09%  '
$ws.Range("E12").Value = '  -2.83%  '
$ws.Range("D13").Value = '1.709.25'
$ws.Range("E13").Value = '  -1.38%  '
$ws.Range("D14").Value = '1.946.11'
$ws.Range("E14").Value = '  -1.18%  '
$ws.Range("E15").Value = '  -2.93%  '
$ws.Range("D16").Value = '0.0₅8192'
$ws.Range("E16").Value = '  -2.35%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.86'
$ws.Range("E17").Value = '  -1.59%  '
$ws.Range("D18").Value = '27.397.64'
$ws.Range("E18").Value = '  -1.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '218.57'
$ws.Range("E19").Value = '  -3.35%  '
$ws.Range("E20").Value = '  +0.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.656'
$ws.Range("E21").Value = '  -3.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.45'
$ws.Range("E22").Value = '  -3.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.956'
$ws.Range("E23").Value = '  -3.96%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.003'
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '142.53'
$ws.Range("E25").Value = '  -3.13%  '
$ws.Range("E26").Value = '  +0.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1211'
$ws.Range("E27").Value = '  -3.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.257'
$ws.Range("E28").Value = '  -2.73%  '
$ws.Range("E29").Value = '  -4.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05404'
$ws.Range("E30").Value = '  -4.58%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.501'
$ws.Range("E32").Value = '  -4.20%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.428'
$ws.Range("E33").Value = '  -2.11%  '
$ws.Range("E34").Value = '  -1.68%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.877'
$ws.Range("E35").Value = '  +1.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9473'
$ws.Range("E36").Value = '  -2.64%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.412'
$ws.Range("E37").Value = '  -0.93%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5878'
$ws.Range("E38").Value = '  -1.65%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01636'
$ws.Range("E39").Value = '  -2.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.854'
$ws.Range("E40").Value = '  -0.91%  '
$ws.Range("D41").Value = '1.048.66'
$ws.Range("E41").Value = '  +0.02%  '
$ws.Range("E42").Value = '  +0.14%  '
$ws.Range("E43").Value = '  -2.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.31'
$ws.Range("E44").Value = '  -0.24%  '
$ws.Range("D45").Value = '1.852.72'
$ws.Range("E45").Value = '  -1.14%  '
$ws.Range("D46").Value = '0.0₈119'
$ws.Range("E46").Value = '  +2.43%  '
$ws.Range("E47").Value = '  -2.65%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4523'
$ws.Range("E48").Value = '  +2.12%  '
$ws.Range("E49").Value = '  +0.40%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.085'
$ws.Range("E50").Value = '  -2.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05226'
$ws.Range("E51").Value = '  -1.98%  '
